$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.1666666666666667
$ws.Range("C2").Value = 0.5930232558139535
$ws.Range("J2").Value = 0.01162790697674419
$ws.Range("P2").Value = 0.1124031007751938
$ws.Range("S2").Value = 0.1162790697674419
$ws.Range("B3").Value = 0.01226993865030675
$ws.Range("C3").Value = 0.04294478527607362
$ws.Range("J3").Value = 0.0245398773006135
$ws.Range("P3").Value = 0.7361963190184049
$ws.Range("S3").Value = 0.1840490797546012
$ws.Range("J4").Value = 0.04444444444444445
$ws.Range("P4").Value = 0.6888888888888889
$ws.Range("S4").Value = 0.2666666666666667
$ws.Range("B6").Value = 0.05607476635514019
$ws.Range("D6").Value = 0.02336448598130841
$ws.Range("F6").Value = 0.06542056074766354
$ws.Range("J6").Value = 0.2570093457943925
$ws.Range("O6").Value = 0.03271028037383177
$ws.Range("Q6").Value = 0.2009345794392523
$ws.Range("R6").Value = 0.0514018691588785
$ws.Range("S6").Value = 0.3130841121495327
$ws.Range("B7").Value = 0.1136363636363636
$ws.Range("D7").Value = 0.01818181818181818
$ws.Range("F7").Value = 0.05
$ws.Range("J7").Value = 0.1590909090909091
$ws.Range("O7").Value = 0.02727272727272727
$ws.Range("Q7").Value = 0.1636363636363636
$ws.Range("S7").Value = 0.3772727272727273
$ws.Range("B8").Value = 0.09504950495049505
$ws.Range("D8").Value = 0.01584158415841584
$ws.Range("F8").Value = 0.06138613861386139
$ws.Range("J8").Value = 0.1148514851485149
$ws.Range("O8").Value = 0.01188118811881188
$ws.Range("Q8").Value = 0.2198019801980198
$ws.Range("R8").Value = 0.08514851485148515
$ws.Range("S8").Value = 0.3960396039603961
$ws.Range("B9").Value = 0.05676855895196507
$ws.Range("F9").Value = 0.06550218340611354
$ws.Range("J9").Value = 0.1135371179039301
$ws.Range("O9").Value = 0.01310043668122271
$ws.Range("Q9").Value = 0.222707423580786
$ws.Range("R9").Value = 0.06550218340611354
$ws.Range("S9").Value = 0.462882096069869
$ws.Range("B10").Value = 0.08890701468189233
$ws.Range("D10").Value = 0.02365415986949429
$ws.Range("E10").Value = 0.003262642740619902
$ws.Range("F10").Value = 0.06933115823817292
$ws.Range("J10").Value = 0.1158238172920065
$ws.Range("O10").Value = 0.01141924959216966
$ws.Range("Q10").Value = 0.2104404567699837
$ws.Range("R10").Value = 0.08319738988580751
$ws.Range("S10").Value = 0.3939641109298532
$ws.Range("G11").Value = 0.1314984709480122
$ws.Range("J11").Value = 0.09480122324159021
$ws.Range("K11").Value = 0.2048929663608563
$ws.Range("L11").Value = 0.5504587155963303
$ws.Range("S11").Value = 0.01834862385321101
$ws.Range("G12").Value = 0.7473118279569892
$ws.Range("J12").Value = 0.1559139784946237
$ws.Range("K12").Value = 0.005376344086021506
$ws.Range("L12").Value = 0.02688172043010753
$ws.Range("S12").Value = 0.06451612903225806
$ws.Range("G13").Value = 0.6911764705882353
$ws.Range("J13").Value = 0.1764705882352941
$ws.Range("S13").Value = 0.1323529411764706
$ws.Range("F15").Value = 0.0186046511627907
$ws.Range("H15").Value = 0.1674418604651163
$ws.Range("I15").Value = 0.09767441860465116
$ws.Range("J15").Value = 0.3209302325581395
$ws.Range("K15").Value = 0.07441860465116279
$ws.Range("M15").Value = 0.02325581395348837
$ws.Range("O15").Value = 0.04651162790697674
$ws.Range("S15").Value = 0.2511627906976744
$ws.Range("F16").Value = 0.01129943502824859
$ws.Range("H16").Value = 0.1864406779661017
$ws.Range("I16").Value = 0.07344632768361582
$ws.Range("J16").Value = 0.4463276836158192
$ws.Range("K16").Value = 0.0903954802259887
$ws.Range("M16").Value = 0.02824858757062147
$ws.Range("N16").Value = 0.005649717514124294
$ws.Range("O16").Value = 0.05084745762711865
$ws.Range("S16").Value = 0.1073446327683616
$ws.Range("F17").Value = 0.01006036217303823
$ws.Range("H17").Value = 0.2173038229376258
$ws.Range("I17").Value = 0.08853118712273642
$ws.Range("J17").Value = 0.358148893360161
$ws.Range("K17").Value = 0.1106639839034205
$ws.Range("M17").Value = 0.03621730382293763
$ws.Range("N17").Value = 0.002012072434607646
$ws.Range("O17").Value = 0.05432595573440644
$ws.Range("S17").Value = 0.1227364185110664
$ws.Range("F18").Value = 0.02094240837696335
$ws.Range("H18").Value = 0.2198952879581152
$ws.Range("I18").Value = 0.07853403141361257
$ws.Range("J18").Value = 0.356020942408377
$ws.Range("K18").Value = 0.1047120418848168
$ws.Range("M18").Value = 0.03664921465968586
$ws.Range("O18").Value = 0.02094240837696335
$ws.Range("S18").Value = 0.162303664921466
$ws.Range("F19").Value = 0.01312910284463895
$ws.Range("H19").Value = 0.2100656455142232
$ws.Range("I19").Value = 0.1013858497447119
$ws.Range("J19").Value = 0.3362509117432531
$ws.Range("K19").Value = 0.1064916119620715
$ws.Range("M19").Value = 0.02479941648431802
$ws.Range("N19").Value = 0.0007293946024799417
$ws.Range("O19").Value = 0.07439824945295405
$ws.Range("S19").Value = 0.1327498176513494

Write-Output "Applied 108 cell updates"